$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.012.86"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "3.893.28"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'468.38"
$ws.Range("E5").Value = "  +9.14%  "
$ws.Range("D6").Value = "'143.32"
$ws.Range("E6").Value = "  +2.58%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.736"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  +6.70%  "
$ws.Range("D11").Value = "'0.0000335"
$ws.Range("E11").Value = "  +6.12%  "
$ws.Range("D12").Value = "'42.82"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "4.512.55"
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").Value = "'10.34"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "'15.20"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "3.872.67"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "'19.89"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("D20").Value = "67.150.04"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "'430.82"
$ws.Range("E21").Value = "  +5.85%  "
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("E23").Value = "  +3.41%  "
$ws.Range("D24").Value = "'88.28"
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("D25").Value = "'38.90"
$ws.Range("E25").Value = "  +5.87%  "
$ws.Range("D26").Value = "'3.53"
$ws.Range("E26").Value = "  +6.71%  "
$ws.Range("E27").Value = "  +6.03%  "
$ws.Range("D28").Value = "'10.04"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "'730.20"
$ws.Range("E30").Value = "  +2.50%  "
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").Value = "'2.79"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value = "'42.90"
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("E35").Value = "  +4.12%  "
$ws.Range("D36").Value = "'57.91"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.37"
$ws.Range("E38").Value = "  -4.55%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0773"
$ws.Range("E39").Value = "  +12.73%  "
$ws.Range("D40").Value = "'3.19"
$ws.Range("E40").Value = "  +10.03%  "
$ws.Range("D41").Value = "'0.0474"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").Value = "'0.338"
$ws.Range("E42").Value = "  +4.29%  "
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "'2.79"
$ws.Range("E45").Value = "  +5.87%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.50"
$ws.Range("E46").Value = "  -6.40%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'2.16"
$ws.Range("E47").Value = "  +3.97%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").Value = "'3.14"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").Value = "'2.88"
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("D51").Value = "'143.50"
$ws.Range("E51").Value = "  +0.74%  "
